$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.625.71'
$ws.Range('E2').Value = '  -1.14%  '
$ws.Range('D3').Value = '1.844.87'
$ws.Range('E3').Value = '  -0.67%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '315.84'
$ws.Range('E5').Value = '  -0.35%  '
$ws.Range('E6').Value = '  -0.06%  '
$ws.Range('E7').Value = '  -2.95%  '
$ws.Range('E8').Value = '  -1.15%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '45.31'
$ws.Range('E9').Value = '  +0.77%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.07260'
$ws.Range('E10').Value = '  -3.01%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.8891'
$ws.Range('E11').Value = '  -5.21%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '20.65'
$ws.Range('E12').Value = '  -3.20%  '
$ws.Range('D13').Value = '1.825.05'
$ws.Range('E13').Value = '  -5.61%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.563'
$ws.Range('E14').Value = '  -2.11%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '5.332'
$ws.Range('E15').Value = '  -1.69%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.06849'
$ws.Range('E16').Value = '  -0.27%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.001'
$ws.Range('E17').Value = '  -0.03%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '78.99'
$ws.Range('E18').Value = '  -2.93%  '
$ws.Range('E19').Value = '  -1.59%  '
$ws.Range('E20').Value = '  -0.13%  '
$ws.Range('E21').Value = '  -2.70%  '
$ws.Range('D22').Value = '27.600.67'
$ws.Range('E22').Value = '  -1.18%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.971'
$ws.Range('E23').Value = '  -2.48%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '10.53'
$ws.Range('E24').Value = '  -4.32%  '
$ws.Range('D25').Value = '2.055.27'
$ws.Range('E25').Value = '  -0.90%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.963'
$ws.Range('E26').Value = '  -2.26%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '155.08'
$ws.Range('E27').Value = '  +0.60%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '18.62'
$ws.Range('E28').Value = '  +1.60%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '123.36'
$ws.Range('E29').Value = '  +9.02%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '5.234'
$ws.Range('E30').Value = '  -3.05%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.853'
$ws.Range('E31').Value = '  +7.12%  '
$ws.Range('E32').Value = '  -0.62%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.7787'
$ws.Range('E33').Value = '  -3.13%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.570'
$ws.Range('E34').Value = '  -5.14%  '
$ws.Range('E35').Value = '  -2.17%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.093'
$ws.Range('E36').Value = '  -6.39%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.9999'
$ws.Range('E37').Value = '  -0.15%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.099'
$ws.Range('E38').Value = '  -1.50%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.05383'
$ws.Range('E39').Value = '  -0.65%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.01928'
$ws.Range('E40').Value = '  -1.97%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.793'
$ws.Range('E41').Value = '  -4.51%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '6.863'
$ws.Range('E42').Value = '  -2.18%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.5066'
$ws.Range('E43').Value = '  -3.28%  '
$ws.Range('E44').Value = '  -1.70%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '8.272'
$ws.Range('E45').Value = '  -5.26%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.06597'
$ws.Range('E46').Value = '  -1.75%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '10.33'
$ws.Range('E47').Value = '  -2.76%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.4712'
$ws.Range('E48').Value = '  -3.11%  '
$ws.Range('E49').Value = '  -2.10%  '
$ws.Range('E50').Value = '  -0.04%  '
$ws.Range('E51').Value = '  -2.58%  '
